$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 759 (shifts the existing rows 759-800
# down to 760-801, same as Excel's UI "Insert Row" / right-click > Insert).
$ws.Rows.Item(759).Insert()

# Populate the newly-inserted row with the new daily-ranking record.
# The date column stores plain text like "2026/12/29" elsewhere in the
# sheet (not a real Excel date), so force text via a leading apostrophe
# and then reset the cell style back to Normal so no stray number-format
# (e.g. date / quoted-text) sticks to the cell.
$ws.Range("A759").Value = "'2026/02/01"
$ws.Range("A759").Style = "Normal"
$ws.Range("B759").Value = "日"
$ws.Range("C759").Value = 3
$ws.Range("D759").Value = 201
